# Upload new version with timestamp
# Two new medicine rows are inserted (alphabetically) into the report table:
#   - ALKAPRESS PLUS 10/160MG 20 F.C. TABS.  (becomes new row 4, before ELICA-M)
#   - URIPAN 5MG 30 TAB.                      (becomes new row 12, before URSOFALK)
# Every other data row shifts down accordingly, the running total is updated,
# and the summary/footer rows move from 17/18 down to 19/20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert two blank rows right before the old "total" row (17),
#    copying the format (styles + row height) from the last data row (16) so
#    the new rows look identical to their neighbours.
$ws.Range("A17:N18").Insert(-4121)

$ws.Range("A16:N16").Copy()
$ws.Range("A17:N17").PasteSpecial(-4122)
$ws.Range("A16:N16").Copy()
$ws.Range("A18:N18").PasteSpecial(-4122)

$ws.Rows(17).RowHeight = 25.5
$ws.Rows(18).RowHeight = 25.5

# 2) Re-create the merged cells for the two new data rows (Insert does not
#    carry merges over from the copied range).
$ws.Range("B17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("B18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()

# 3) Rewrite the full data block (rows 4-18) with the final, sorted contents:
#    serial number (A), item name (B), transactions summary (H), balance (L)
#    and price summary (N).
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ALKAPRESS PLUS 10/160MG 20 F.C. TABS."
$ws.Range("H4").Value = "0:1"
$ws.Range("L4").Value = 51
$ws.Range("N4").Value = "0:2"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "ELICA-M CREAM 30 GRAM"
$ws.Range("H5").Value = "0:0"
$ws.Range("L5").Value = 52
$ws.Range("N5").Value = "1:0"

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "KERELLA LOTION 30 ML"
$ws.Range("H6").Value = "3:0"
$ws.Range("L6").Value = 31
$ws.Range("N6").Value = "1:0"

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "LIDOCAINE 10% TOPICAL SPRAY 15 GM"
$ws.Range("H7").Value = "1:0"
$ws.Range("L7").Value = 38
$ws.Range("N7").Value = "1:0"

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "NEUROGLOPENTIN 300 MG 30 CAPS."
$ws.Range("H8").Value = "1:0"
$ws.Range("L8").Value = 37
$ws.Range("N8").Value = "0:0"

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "SULBIN 750MG VIAL"
$ws.Range("H9").Value = "6:0"
$ws.Range("L9").Value = 35
$ws.Range("N9").Value = "1:0"

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "SUPOLACK HAIR SHAMPOO 200 ML"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 149.5
$ws.Range("N10").Value = "1:0"

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "TELFAST 180MG 20 F.C. TABS"
$ws.Range("H11").Value = "1:0"
$ws.Range("L11").Value = 80
$ws.Range("N11").Value = "0:2"

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "URIPAN 5MG 30 TAB."
$ws.Range("H12").Value = "1:1"
$ws.Range("L12").Value = 18
$ws.Range("N12").Value = "0:0"

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "URSOFALK 250MG 20 CAPS."
$ws.Range("H13").Value = "0:0"
$ws.Range("L13").Value = 122
$ws.Range("N13").Value = "1:0"

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
$ws.Range("H14").Value = "6:0"
$ws.Range("L14").Value = 26
$ws.Range("N14").Value = "1:0"

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "WELLMETAZONE 0.1% CREAM 40 GM"
$ws.Range("H15").Value = "0:0"
$ws.Range("L15").Value = 56
$ws.Range("N15").Value = "1:0"

$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "جنتيانا نقط"
$ws.Range("H16").Value = "4:0"
$ws.Range("L16").Value = 14
$ws.Range("N16").Value = "2:0"

$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "سرنجات 5 سم"
$ws.Range("H17").Value = "-1:0"
$ws.Range("L17").Value = 2
$ws.Range("N17").Value = "1:0"

$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "كريم فاتيكا 125 مل"
$ws.Range("H18").Value = "2:0"
$ws.Range("L18").Value = 50
$ws.Range("N18").Value = "1:0"

# 4) Update the running total (old row 17, now row 19) to reflect the two
#    newly added rows (692.5 + 51 + 18 = 761.5).
$ws.Range("K19").Value = 761.5
